# Apr 1 data update for sd_covid_dataset
#
# 1. A new "paumavalley" column is inserted just before the existing
#    "ramona" column (old column AZ / 52), shifting ramona, ranchosantafe,
#    sanysidro, springvalley and other one column to the right.
# 2. tested (col B) gets a value for row 27 (Mar 31).
# 3. A brand new row 28 (Apr 1, 2020) is appended with that day's counts.
# 4. The hidden _FilterDatabase name is widened to cover the new column.
# 5. Selection is moved to the last-entered cell, matching a normal
#    "finished typing the new row" interactive session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new "paumavalley" column before old column AZ (52) ---
$ws.Columns("AZ:AZ").Insert()
$ws.Range("AZ1").Value = "paumavalley"

# --- 2. Fill in the missing "tested" figure for Mar 31 (row 27) ---
$ws.Range("B27").Value = 989

# --- 3. Append the Apr 1, 2020 row (row 28) ---
$row = 28
$ws.Cells.Item($row, 1).Value = 43922   # date -> 2020-04-01

$rowValues = @{
    3  = 849
    4  = 6
    5  = 0
    6  = 8
    7  = 1
    8  = 163
    9  = 10
    10 = 194
    11 = 18
    12 = 154
    13 = 19
    14 = 125
    15 = 31
    16 = 93
    17 = 27
    18 = 62
    19 = 27
    20 = 41
    21 = 25
    22 = 3
    23 = 0
    24 = 380
    25 = 461
    26 = 8
    27 = 158
    28 = 65
    29 = 15
    30 = 30
    31 = 52
    32 = 2
    33 = 5
    34 = 38
    35 = 24
    36 = 14
    38 = 11
    39 = 5
    40 = 11
    41 = 20
    42 = 9
    43 = 486
    44 = 10
    45 = 5
    46 = 3
    47 = 9
    48 = 6
    49 = 4
    51 = 8
    52 = 1
    53 = 4
    54 = 12
    56 = 18
    57 = 34
}

foreach ($col in $rowValues.Keys) {
    $ws.Cells.Item($row, $col).Value = $rowValues[$col]
}

# --- 4. Widen the hidden AutoFilter-database name to include the new column ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$BB`$17"
    }
}

# --- 5. Leave the selection on the last cell touched, like a live edit session ---
$ws.Range("Y28").Select()
